$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value2 = 1554.4546
$ws.Range("I28").Value2 = 2181.5
$ws.Range("K28").Value2 = 2181.5
$ws.Range("M28").Value2 = -1696.5

$ws.Range("H98").Value2 = 2603.8
$ws.Range("I98").Value2 = 2337.5557
$ws.Range("J98").Value2 = 5000
$ws.Range("K98").Value2 = 2337.5557
$ws.Range("L98").Value2 = 5000
$ws.Range("M98").Value2 = -839.5556999999999
$ws.Range("N98").Value2 = -7996

$ws.Range("H122").Value2 = 2603.8
$ws.Range("I122").Value2 = 2337.5557
$ws.Range("J122").Value2 = 5000
$ws.Range("K122").Value2 = 7012.6671
$ws.Range("L122").Value2 = 15000
$ws.Range("M122").Value2 = -4562.6671
$ws.Range("N122").Value2 = -19900

$ws.Range("H127").Value2 = 1174.5
$ws.Range("I127").Value2 = 200
$ws.Range("J127").Value2 = 1499.3334
$ws.Range("K127").Value2 = 600
$ws.Range("L127").Value2 = 4498.0002
$ws.Range("M127").Value2 = 4360
$ws.Range("N127").Value2 = -14418.0002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value2 = 1649
$ws.Range("I2").Value2 = 1690.6
$ws.Range("J2").Value2 = 1545
$ws.Range("K2").Value2 = 1690.6
$ws.Range("L2").Value2 = 1545
$ws.Range("M2").Value2 = -1577.6
$ws.Range("N2").Value2 = -1771

$ws.Range("H61").Value2 = 3012
$ws.Range("I61").Value2 = 3012
$ws.Range("J61").Value2 = 0
$ws.Range("K61").Value2 = 3012
$ws.Range("L61").Value2 = 0
$ws.Range("M61").ClearContents() | Out-Null
$ws.Range("N61").Value2 = -2800

$ws.Range("H62").Value2 = 9750
$ws.Range("J62").Value2 = 9750
$ws.Range("L62").Value2 = 9750
$ws.Range("N62").Value2 = -10998

$ws.Range("H65").Value2 = 9750
$ws.Range("J65").Value2 = 9750
$ws.Range("L65").Value2 = 29250
$ws.Range("N65").Value2 = -35490

$ws.Range("H94").Value2 = 39552.168
$ws.Range("J94").Value2 = 39552.168
$ws.Range("L94").Value2 = 39552.168
$ws.Range("N94").Value2 = -41354.168

$ws.Range("H101").Value2 = 0
$ws.Range("J101").Value2 = 0
$ws.Range("L101").ClearContents() | Out-Null
$ws.Range("N101").Value2 = 0

$ws.Range("H102").Value2 = 2995.6
$ws.Range("I102").Value2 = 2995.6
$ws.Range("K102").Value2 = 2995.6
$ws.Range("M102").Value2 = -1373.6

$ws.Range("H107").Value2 = 0
$ws.Range("J107").Value2 = 0
$ws.Range("L107").ClearContents() | Out-Null
$ws.Range("N107").Value2 = 0

$ws.Range("H116").Value2 = 1649
$ws.Range("I116").Value2 = 1690.6
$ws.Range("J116").Value2 = 1545
$ws.Range("K116").Value2 = 1690.6
$ws.Range("L116").Value2 = 1545
$ws.Range("M116").Value2 = 603.4000000000001
$ws.Range("N116").Value2 = -6133

$ws.Range("H136").Value2 = 3012
$ws.Range("I136").Value2 = 3012
$ws.Range("J136").Value2 = 0
$ws.Range("K136").Value2 = 9036
$ws.Range("L136").Value2 = 0
$ws.Range("M136").ClearContents() | Out-Null
$ws.Range("N136").Value2 = -6486

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value2 = 1649
$ws.Range("I3").Value2 = 1690.6
$ws.Range("J3").Value2 = 1545
$ws.Range("K3").Value2 = 1690.6
$ws.Range("L3").Value2 = 1545
$ws.Range("M3").Value2 = -1576.6
$ws.Range("N3").Value2 = -1773

$ws.Range("H64").Value2 = 3715.6
$ws.Range("I64").Value2 = 344
$ws.Range("J64").Value2 = 4558.5
$ws.Range("K64").Value2 = 344
$ws.Range("L64").Value2 = 4558.5
$ws.Range("M64").Value2 = -119
$ws.Range("N64").Value2 = -5008.5

$ws.Range("H67").Value2 = 3715.6
$ws.Range("I67").Value2 = 344
$ws.Range("J67").Value2 = 4558.5
$ws.Range("K67").Value2 = 344
$ws.Range("L67").Value2 = 4558.5
$ws.Range("M67").Value2 = 436
$ws.Range("N67").Value2 = -6118.5

$ws.Range("H99").Value2 = 2531.6
$ws.Range("I99").Value2 = 2164.5
$ws.Range("K99").Value2 = 2164.5
$ws.Range("M99").Value2 = -666.5

$ws.Range("H100").Value2 = 36000
$ws.Range("J100").Value2 = 36000
$ws.Range("L100").Value2 = 36000
$ws.Range("N100").Value2 = -38164

$ws.Range("H107").Value2 = 2227.5
$ws.Range("J107").Value2 = 2999.5
$ws.Range("L107").Value2 = 2999.5
$ws.Range("N107").Value2 = -6839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value2 = 72249
$ws.Range("J53").Value2 = 72249
$ws.Range("L53").Value2 = 72249
$ws.Range("N53").Value2 = -73463

$ws.Range("H88").Value2 = 9960
$ws.Range("J88").Value2 = 9960
$ws.Range("L88").Value2 = 9960
$ws.Range("N88").Value2 = -10772

$ws.Range("H91").Value2 = 9960
$ws.Range("J91").Value2 = 9960
$ws.Range("L91").Value2 = 9960
$ws.Range("N91").Value2 = -12768

$ws.Range("H122").Value2 = 1700
$ws.Range("I122").Value2 = 1700
$ws.Range("K122").Value2 = 5100
$ws.Range("M122").Value2 = -2650

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value2 = 706
$ws.Range("I5").Value2 = 483.2857
$ws.Range("J5").Value2 = 1017.8
$ws.Range("K5").Value2 = 1449.8571
$ws.Range("L5").Value2 = 3053.4
$ws.Range("M5").Value2 = -1337.8571
$ws.Range("N5").Value2 = -3277.4

$ws.Range("H11").Value2 = 50.125
$ws.Range("I11").Value2 = 47
$ws.Range("J11").Value2 = 53.25
$ws.Range("K11").Value2 = 141
$ws.Range("L11").Value2 = 159.75
$ws.Range("M11").Value2 = -1
$ws.Range("N11").Value2 = -439.75

$ws.Range("H12").Value2 = 218.8421
$ws.Range("I12").Value2 = 392.22223
$ws.Range("J12").Value2 = 62.8
$ws.Range("K12").Value2 = 1176.66669
$ws.Range("L12").Value2 = 188.4
$ws.Range("M12").Value2 = -1003.66669
$ws.Range("N12").Value2 = -534.4

$ws.Range("H14").Value2 = 429.33334
$ws.Range("I14").Value2 = 429.33334
$ws.Range("K14").Value2 = 1288.00002
$ws.Range("M14").Value2 = -1115.00002

$ws.Range("H62").Value2 = 0
$ws.Range("J62").Value2 = 0
$ws.Range("L62").ClearContents() | Out-Null
$ws.Range("N62").Value2 = 0

$ws.Range("H65").Value2 = 0
$ws.Range("J65").Value2 = 0
$ws.Range("L65").ClearContents() | Out-Null
$ws.Range("N65").Value2 = 0

$ws.Range("H68").Value2 = 998.75
$ws.Range("I68").Value2 = 998
$ws.Range("K68").Value2 = 2994
$ws.Range("M68").Value2 = -2183

$ws.Range("H71").Value2 = 998.75
$ws.Range("I71").Value2 = 998
$ws.Range("K71").Value2 = 8982
$ws.Range("M71").Value2 = -4926

$ws.Range("H100").Value2 = 4200
$ws.Range("J100").Value2 = 0
$ws.Range("L100").Value2 = 0
$ws.Range("N100").ClearContents() | Out-Null

$ws.Range("H113").Value2 = 1542.5294
$ws.Range("I113").Value2 = 957.2
$ws.Range("J113").Value2 = 1786.4166
$ws.Range("K113").Value2 = 2871.6
$ws.Range("L113").Value2 = 5359.2498
$ws.Range("M113").Value2 = -701.6000000000004
$ws.Range("N113").Value2 = -9699.2498

$ws.Range("H120").Value2 = 999
$ws.Range("I120").Value2 = 999
$ws.Range("K120").Value2 = 2997
$ws.Range("M120").Value2 = 1841

$ws.Range("H132").Value2 = 2495.6667
$ws.Range("I132").Value2 = 1500
$ws.Range("J132").Value2 = 2993.5
$ws.Range("K132").Value2 = 13500
$ws.Range("L132").Value2 = 26941.5
$ws.Range("M132").Value2 = -10970
$ws.Range("N132").Value2 = -32001.5

$ws.Range("H135").Value2 = 706
$ws.Range("I135").Value2 = 483.2857
$ws.Range("J135").Value2 = 1017.8
$ws.Range("K135").Value2 = 4349.571300000001
$ws.Range("L135").Value2 = 9160.199999999999
$ws.Range("M135").Value2 = -1814.571300000001
$ws.Range("N135").Value2 = -14230.2

$ws.Range("H141").Value2 = 32682.5
$ws.Range("I141").Value2 = 32682.5
$ws.Range("K141").Value2 = 98047.5
$ws.Range("M141").Value2 = -92867.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value2 = 8500
$ws.Range("I70").Value2 = 8500
$ws.Range("J70").Value2 = 0
$ws.Range("K70").Value2 = 8500
$ws.Range("L70").Value2 = 0
$ws.Range("M70").ClearContents() | Out-Null
$ws.Range("N70").Value2 = -8230

$ws.Range("H73").Value2 = 8500
$ws.Range("I73").Value2 = 8500
$ws.Range("J73").Value2 = 0
$ws.Range("K73").Value2 = 8500
$ws.Range("L73").Value2 = 0
$ws.Range("M73").ClearContents() | Out-Null
$ws.Range("N73").Value2 = -7564

$ws.Range("H97").Value2 = 1099.6666
$ws.Range("I97").Value2 = 332.66666
$ws.Range("J97").Value2 = 1866.6666
$ws.Range("K97").Value2 = 332.66666
$ws.Range("L97").Value2 = 1866.6666
$ws.Range("M97").Value2 = 163.33334
$ws.Range("N97").Value2 = -2858.6666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 7402
$ws.Range("I22").Value2 = 8492.888999999999
$ws.Range("K22").Value2 = 8492.888999999999
$ws.Range("M22").Value2 = -8197.888999999999

$ws.Range("H27").Value2 = 7402
$ws.Range("I27").Value2 = 8492.888999999999
$ws.Range("K27").Value2 = 8492.888999999999
$ws.Range("M27").Value2 = -8385.888999999999

$ws.Range("H68").Value2 = 1296
$ws.Range("I68").Value2 = 1296
$ws.Range("K68").Value2 = 1296
$ws.Range("M68").Value2 = -547

$ws.Range("H71").Value2 = 1296
$ws.Range("I71").Value2 = 1296
$ws.Range("K71").Value2 = 6480
$ws.Range("M71").Value2 = -2736

$ws.Range("H132").Value2 = 4528.8
$ws.Range("J132").Value2 = 5499.5
$ws.Range("L132").Value2 = 16498.5
$ws.Range("N132").Value2 = -21558.5

$ws.Range("H136").Value2 = 1253.8572
$ws.Range("I136").Value2 = 1355.4
$ws.Range("K136").Value2 = 4066.2
$ws.Range("M136").Value2 = -1516.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H93").Value2 = 129999
$ws.Range("I93").Value2 = 129999
$ws.Range("K93").Value2 = 129999
$ws.Range("M93").Value2 = -127503

$ws.Range("H100").Value2 = 1467.8
$ws.Range("I100").Value2 = 1899.5
$ws.Range("J100").Value2 = 1180
$ws.Range("K100").Value2 = 3799
$ws.Range("L100").Value2 = 2360
$ws.Range("M100").Value2 = -3258
$ws.Range("N100").Value2 = -3442

$ws.Range("H113").Value2 = 544.6667
$ws.Range("I113").Value2 = 533.6
$ws.Range("K113").Value2 = 1600.8
$ws.Range("M113").Value2 = 569.1999999999998
